$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (B2:G2, I2:K2)
$ws.Range("B2").Value = 4.938
$ws.Range("C2").Value = 66.63
$ws.Range("D2").Value = 59.243
$ws.Range("E2").Value = -80.269
$ws.Range("F2").Value = 62.379
$ws.Range("G2").Value = 77.441
$ws.Range("I2").Value = 3.369
$ws.Range("J2").Value = 0.054
$ws.Range("K2").Value = 7.825

# H2 keeps using a shared string, but the shared string content itself changes.
$ws.Range("H2").Value = "[   0.7238   -0.0505   -0.278  -122.4354   42.8765  179.7921]"

# Remove row 3 entirely
$ws.Rows("3:3").Delete()
